$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 115, shifting the existing rows 115:123 down to 116:124
$ws.Rows.Item(115).Insert()

# Populate the newly inserted row 115 with the new weekly data, copying the
# constant columns from the row above (row 114) and setting the varying ones.
$ws.Range("A115").Value = 7
$ws.Range("B115").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C115").Value = "Ñuble"
$ws.Range("D115").Value = 45041
$ws.Range("D115").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E115").Value = 16
$ws.Range("F115").Value = "Fruta"
$ws.Range("G115").Value = 100108
$ws.Range("H115").Value = "Tropicales y subtropicales"
$ws.Range("I115").Value = 100108002
$ws.Range("J115").Value = "Mango"
$ws.Range("K115").Value = "Sin especificar"
$ws.Range("L115").Value = "Primera"
$ws.Range("M115").Value = 90
$ws.Range("N115").Value = 7000
$ws.Range("O115").Value = 8000
$ws.Range("P115").Value = 7556
$ws.Range("Q115").Value = "$/bandeja 4 kilos"
$ws.Range("R115").Value = "Perú"
$ws.Range("S115").Value = 1889
$ws.Range("T115").Value = 4
